# Update benchmark result values on Sheet1 (data for rows 2-5, columns B-G)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 110.0884997844696
$ws.Range("C2").Value = 0.540052929905726
$ws.Range("D2").Value = 0.540052929905726
$ws.Range("E2").Value = 95.08698892593384
$ws.Range("F2").Value = 0.6252550167914027
$ws.Range("G2").Value = 0.6252550167914027

$ws.Range("B3").Value = 43.52274942398071
$ws.Range("C3").Value = 1.366035410087631
$ws.Range("D3").Value = 0.6830177050438153
$ws.Range("E3").Value = 41.95485782623291
$ws.Range("F3").Value = 1.417085408888084
$ws.Range("G3").Value = 0.7085427044440419

$ws.Range("B4").Value = 25.15734624862671
$ws.Range("C4").Value = 2.363270603741606
$ws.Range("D4").Value = 0.5908176509354016
$ws.Range("E4").Value = 25.50388050079346
$ws.Range("F4").Value = 2.331159638850998
$ws.Range("G4").Value = 0.5827899097127495

$ws.Range("B5").Value = 21.81533575057983
$ws.Range("C5").Value = 2.725312942109931
$ws.Range("D5").Value = 0.3406641177637414
$ws.Range("E5").Value = 24.43703198432922
$ws.Range("F5").Value = 2.432931171660069
$ws.Range("G5").Value = 0.3041163964575087
